$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) B1:E1
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2) B2:E2
$ws.Range("B2").Value = 15.770593801060102
$ws.Range("C2").Value = 20.493003485450572
$ws.Range("D2").Value = 28.479259450375366
$ws.Range("E2").Value = 30.801006908779186

# Update STR row (row 3) B3:E3
$ws.Range("B3").Value = 13.102229730301303
$ws.Range("C3").Value = 22.44654022273221
$ws.Range("D3").Value = 19.578343119659678
$ws.Range("E3").Value = 33.33869612715057

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
